# Daily attendance processing - 2025-12-25 12:50:35
# Reorders the "Recorded By" (column G) contributor lists for several rows
# on the active worksheet. The set of recorders per row is unchanged; only
# the order in which they are listed is updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2"  = "servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, System, Veronia.rafat@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    "G3"  = "eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, System, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg"
    "G4"  = "servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg"
    "G5"  = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
    "G6"  = "alshimaa.atef@med.asu.edu.egm, Mohammedeltanany@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, manar.montaser@med.asu.edu.eg"
    "G7"  = "AbeerRagheb@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg"
    "G12" = "yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"
    "G13" = "amira.m.ibrahim@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, esraa.mostafa@med.asu.edu.eg"
    "G24" = "youstina.gamil@med.asu.edu.eg, Sarah.Mahdy@med.asu.edu.eg"
    "G25" = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"
    "G28" = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
    "G30" = "wafaa.ebida@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
